$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("C3").Value = 155732
$ws.Range("C4").Value = 146836
$ws.Range("C5").Value = 8896
$ws.Range("C8").Value = 63.73
